# Update the "Wendy - Socioeconomic Disparities in Air Quality and Health" bullet
# on slide 1 with the rewritten hypothesis statement.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

$target = $null
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $shape = $s.Shapes.Item($i)
    if ($shape.HasTextFrame -and $shape.TextFrame.HasText) {
        if ($shape.TextFrame.TextRange.Text -like "*Socioeconomic Disparities in Air Quality and Health*") {
            $target = $shape
            break
        }
    }
}

$tr = $target.TextFrame.TextRange
$oldText = ": Explore whether there are disparities in air quality and health outcomes across different socioeconomic groups. This could involve analyzing data on income levels, education levels, and race/ethnicity alongside air quality and health data to identify any correlations or disparities."
$newText = ": Socioeconomic status, as indicated by income levels, education attainment, and race/ethnicity, is a significant predictor of air quality and health outcomes. Communities with lower socioeconomic status are hypothesized to experience poorer air quality, which in turn leads to a higher prevalence of adverse health outcomes. This relationship is expected to persist even when controlling for potential confounding variables such as geographic location and access to healthcare services."

$paraCount = $tr.Paragraphs().Count
for ($pi = 1; $pi -le $paraCount; $pi++) {
    $para = $tr.Paragraphs($pi, 1)
    $runCount = $para.Runs().Count
    for ($ri = 1; $ri -le $runCount; $ri++) {
        $run = $para.Runs($ri, 1)
        if ($run.Text -eq $oldText) {
            $run.Text = $newText
        }
    }
}
